$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 57.045267
$ws.Range("H2").Value = 171.135801
$ws.Range("I2").Value = 0.2489699905037019
$ws.Range("J2").Value = 0.2489699905037019
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 1.140304
$ws.Range("N2").Value = 3.420912
$ws.Range("O2").Value = 0.263324902255146
$ws.Range("P2").Value = 0.263324902255146
$ws.Range("Q2").Value = 65.04894614116802
$ws.Range("R2").Value = 585.4405152705122
$ws.Range("S2").Value = 0.06555999841385193
$ws.Range("T2").Value = 0.06555999841385192
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 57.045267
$ws.Range("H3").Value = 171.135801
$ws.Range("I3").Value = 0.2489699905037019
$ws.Range("J3").Value = 0.2489699905037019
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 1.916907333333333
$ws.Range("N3").Value = 5.750722
$ws.Range("O3").Value = 0.4426621639336287
$ws.Range("P3").Value = 0.4426621639336287
$ws.Range("Q3").Value = 109.350490644258
$ws.Range("R3").Value = 984.154415798322
$ws.Range("S3").Value = 0.1102095947509037
$ws.Range("T3").Value = 0.1102095947509037
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 57.045267
$ws.Range("H4").Value = 171.135801
$ws.Range("I4").Value = 0.2489699905037019
$ws.Range("J4").Value = 0.2489699905037019
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.3556303333333333
$ws.Range("N4").Value = 1.066891
$ws.Range("O4").Value = 0.08212399742872513
$ws.Range("P4").Value = 0.08212399742872513
$ws.Range("Q4").Value = 20.287027318299
$ws.Range("R4").Value = 182.583245864691
$ws.Range("S4").Value = 0.02044641085995573
$ws.Range("T4").Value = 0.02044641085995573
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 57.045267
$ws.Range("H5").Value = 171.135801
$ws.Range("I5").Value = 0.2489699905037019
$ws.Range("J5").Value = 0.2489699905037019
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 0.9175653333333335
$ws.Range("N5").Value = 2.752696
$ws.Range("O5").Value = 0.2118889363825002
$ws.Range("P5").Value = 0.2118889363825001
$ws.Range("Q5").Value = 52.34275942994401
$ws.Range("R5").Value = 471.0848348694961
$ws.Range("S5").Value = 0.05275398647899056
$ws.Range("T5").Value = 0.05275398647899055
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 54.09018966666667
$ws.Range("H6").Value = 162.270569
$ws.Range("I6").Value = 0.2360727666969011
$ws.Range("J6").Value = 0.2360727666969011
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 1.140304
$ws.Range("N6").Value = 3.420912
$ws.Range("O6").Value = 0.263324902255146
$ws.Range("P6").Value = 0.263324902255146
$ws.Range("Q6").Value = 61.67925963765868
$ws.Range("R6").Value = 555.1133367389281
$ws.Range("S6").Value = 0.06216383821556338
$ws.Range("T6").Value = 0.06216383821556336
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 54.09018966666667
$ws.Range("H7").Value = 162.270569
$ws.Range("I7").Value = 0.2360727666969011
$ws.Range("J7").Value = 0.2360727666969011
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 1.916907333333333
$ws.Range("N7").Value = 5.750722
$ws.Range("O7").Value = 0.4426621639336287
$ws.Range("P7").Value = 0.4426621639336287
$ws.Range("Q7").Value = 103.6858812334242
$ws.Range("R7").Value = 933.1729311008179
$ws.Range("S7").Value = 0.1045004817518489
$ws.Range("T7").Value = 0.1045004817518489
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 54.09018966666667
$ws.Range("H8").Value = 162.270569
$ws.Range("I8").Value = 0.2360727666969011
$ws.Range("J8").Value = 0.2360727666969011
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 0.3556303333333333
$ws.Range("N8").Value = 1.066891
$ws.Range("O8").Value = 0.08212399742872513
$ws.Range("P8").Value = 0.08212399742872513
$ws.Range("Q8").Value = 19.23611218121989
$ws.Range("R8").Value = 173.125009630979
$ws.Range("S8").Value = 0.01938723928520833
$ws.Range("T8").Value = 0.01938723928520833
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 54.09018966666667
$ws.Range("H9").Value = 162.270569
$ws.Range("I9").Value = 0.2360727666969011
$ws.Range("J9").Value = 0.2360727666969011
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 0.9175653333333335
$ws.Range("N9").Value = 2.752696
$ws.Range("O9").Value = 0.2118889363825002
$ws.Range("P9").Value = 0.2118889363825001
$ws.Range("Q9").Value = 49.63128291155823
$ws.Range("R9").Value = 446.681546204024
$ws.Range("S9").Value = 0.05002120744428049
$ws.Range("T9").Value = 0.05002120744428047
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 67.444722
$ws.Range("H10").Value = 202.334166
$ws.Range("I10").Value = 0.2943576685488177
$ws.Range("J10").Value = 0.2943576685488177
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 1.140304
$ws.Range("N10").Value = 3.420912
$ws.Range("O10").Value = 0.263324902255146
$ws.Range("P10").Value = 0.263324902255146
$ws.Range("Q10").Value = 76.90748627548801
$ws.Range("R10").Value = 692.167376479392
$ws.Range("S10").Value = 0.0775117042986701
$ws.Range("T10").Value = 0.07751170429867009
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 67.444722
$ws.Range("H11").Value = 202.334166
$ws.Range("I11").Value = 0.2943576685488177
$ws.Range("J11").Value = 0.2943576685488177
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 1.916907333333333
$ws.Range("N11").Value = 5.750722
$ws.Range("O11").Value = 0.4426621639336287
$ws.Range("P11").Value = 0.4426621639336287
$ws.Range("Q11").Value = 129.285282196428
$ws.Range("R11").Value = 1163.567539767852
$ws.Range("S11").Value = 0.1303010025302775
$ws.Range("T11").Value = 0.1303010025302775
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 67.444722
$ws.Range("H12").Value = 202.334166
$ws.Range("I12").Value = 0.2943576685488177
$ws.Range("J12").Value = 0.2943576685488177
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 0.3556303333333333
$ws.Range("N12").Value = 1.066891
$ws.Range("O12").Value = 0.08212399742872513
$ws.Range("P12").Value = 0.08212399742872513
$ws.Range("Q12").Value = 23.985388966434
$ws.Range("R12").Value = 215.868500697906
$ws.Range("S12").Value = 0.02417382841502863
$ws.Range("T12").Value = 0.02417382841502863
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 67.444722
$ws.Range("H13").Value = 202.334166
$ws.Range("I13").Value = 0.2943576685488177
$ws.Range("J13").Value = 0.2943576685488177
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 0.9175653333333335
$ws.Range("N13").Value = 2.752696
$ws.Range("O13").Value = 0.2118889363825002
$ws.Range("P13").Value = 0.2118889363825001
$ws.Range("Q13").Value = 61.88493882350401
$ws.Range("R13").Value = 556.964449411536
$ws.Range("S13").Value = 0.06237113330484151
$ws.Range("T13").Value = 0.0623711333048415
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 50.54489333333333
$ws.Range("H14").Value = 151.63468
$ws.Range("I14").Value = 0.2205995742505793
$ws.Range("J14").Value = 0.2205995742505793
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 1.140304
$ws.Range("N14").Value = 3.420912
$ws.Range("O14").Value = 0.263324902255146
$ws.Range("P14").Value = 0.263324902255146
$ws.Range("Q14").Value = 57.63654404757335
$ws.Range("R14").Value = 518.7288964281601
$ws.Range("S14").Value = 0.05808936132706063
$ws.Range("T14").Value = 0.05808936132706062
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 50.54489333333333
$ws.Range("H15").Value = 151.63468
$ws.Range("I15").Value = 0.2205995742505793
$ws.Range("J15").Value = 0.2205995742505793
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 1.916907333333333
$ws.Range("N15").Value = 5.750722
$ws.Range("O15").Value = 0.4426621639336287
$ws.Range("P15").Value = 0.4426621639336287
$ws.Range("Q15").Value = 96.88987669321777
$ws.Range("R15").Value = 872.00889023896
$ws.Range("S15").Value = 0.09765108490059864
$ws.Range("T15").Value = 0.09765108490059864
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 50.54489333333333
$ws.Range("H16").Value = 151.63468
$ws.Range("I16").Value = 0.2205995742505793
$ws.Range("J16").Value = 0.2205995742505793
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 0.3556303333333333
$ws.Range("N16").Value = 1.066891
$ws.Range("O16").Value = 0.08212399742872513
$ws.Range("P16").Value = 0.08212399742872513
$ws.Range("Q16").Value = 17.97529726443111
$ws.Range("R16").Value = 161.77767537988
$ws.Range("S16").Value = 0.01811651886853244
$ws.Range("T16").Value = 0.01811651886853244
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 50.54489333333333
$ws.Range("H17").Value = 151.63468
$ws.Range("I17").Value = 0.2205995742505793
$ws.Range("J17").Value = 0.2205995742505793
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 0.9175653333333335
$ws.Range("N17").Value = 2.752696
$ws.Range("O17").Value = 0.2118889363825002
$ws.Range("P17").Value = 0.2118889363825001
$ws.Range("Q17").Value = 46.37824189969778
$ws.Range("R17").Value = 417.4041770972801
$ws.Range("S17").Value = 0.04674260915438763
$ws.Range("T17").Value = 0.04674260915438762
